# New datasets + baseline regression
#
# Across this workbook's "variables_####" sheets, column A holds a short
# list of category labels. One label - "congenital" - is being renamed to
# "misc_long_term" everywhere it appears (it shows up in either A3 or A4
# depending on the sheet's layout), while all the other labels are left
# untouched.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    foreach ($r in 1..10) {
        $cell = $ws.Cells.Item($r, 1)
        if ($cell.Text -eq "congenital") {
            $cell.Value = "misc_long_term"
        }
    }
}
